# update ChooseItemManager add Weapon to button
# Adds a "ScriptName" column (H) to the WeaponData sheet, fills it in for
# the existing MagicBall row, and duplicates that row twice (rows 5 & 6)
# so the sheet has three MagicBall entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WeaponData")

# --- New column H header + type/localization rows -------------------------
# Copy the style from column G (same row) onto column H first, then set the
# text values so the new cells pick up the existing formatting.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "ScriptName"

$ws.Range("G2").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("H2").Value = "程式名稱"

$ws.Range("G3").Copy()
$ws.Range("H3").PasteSpecial(-4122)
$ws.Range("H3").Value = "string"

$ws.Range("G4").Copy()
$ws.Range("H4").PasteSpecial(-4122)
$ws.Range("H4").Value = "MagicBallController"

# --- Duplicate the MagicBall data row (row 4) into rows 5 and 6 -----------
$ws.Range("A4:H4").Copy()
$ws.Range("A5").PasteSpecial(-4163)
$ws.Range("A4:H4").Copy()
$ws.Range("A5").PasteSpecial(-4122)

$ws.Range("A4:H4").Copy()
$ws.Range("A6").PasteSpecial(-4163)
$ws.Range("A4:H4").Copy()
$ws.Range("A6").PasteSpecial(-4122)

# --- Column width for the new column ---------------------------------------
$ws.Columns.Item(8).ColumnWidth = 17.417
